$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (year header) ---
$ws.Range("O3").Copy($ws.Range("P3"))
$ws.Range("P3").Value2 = 2022

# --- Row 5 (totals row, bold-ish numeric style like O5) ---
$ws.Range("O5").Copy($ws.Range("P5"))
$ws.Range("P5").Value2 = 2349.7

# --- Rows 6-14 (plain numeric style) ---
$ws.Range("O6").Copy($ws.Range("P6"))
$ws.Range("P6").Value2 = 103.2

$ws.Range("O7").Copy($ws.Range("P7"))
$ws.Range("P7").Value2 = 231.7

$ws.Range("O8").Copy($ws.Range("P8"))
$ws.Range("P8").Value2 = 211.3

$ws.Range("O9").Copy($ws.Range("P9"))
$ws.Range("P9").Value2 = 226.9

$ws.Range("O10").Copy($ws.Range("P10"))
$ws.Range("P10").Value2 = 297

$ws.Range("O11").Copy($ws.Range("P11"))
$ws.Range("P11").Value2 = 321.1

$ws.Range("O12").Copy($ws.Range("P12"))
$ws.Range("P12").Value2 = 923

$ws.Range("O13").Copy($ws.Range("P13"))
$ws.Range("P13").Value2 = 35.4

$ws.Range("O14").Copy($ws.Range("P14"))
$ws.Range("P14").Value2 = "-"

# --- Row 15 (blank section header row) ---
$ws.Range("O15").Copy($ws.Range("P15"))

# --- Rows 16-24 (percentage rows) ---
$ws.Range("O16").Copy($ws.Range("P16"))
$ws.Range("P16").Value2 = 26.7

$ws.Range("O17").Copy($ws.Range("P17"))
$ws.Range("P17").Value2 = 15.1

$ws.Range("O18").Copy($ws.Range("P18"))
$ws.Range("P18").Value2 = 21

$ws.Range("O19").Copy($ws.Range("P19"))
$ws.Range("P19").Value2 = 29.6

$ws.Range("O20").Copy($ws.Range("P20"))
$ws.Range("P20").Value2 = 29.9

$ws.Range("O21").Copy($ws.Range("P21"))
$ws.Range("P21").Value2 = 21.5

$ws.Range("O22").Copy($ws.Range("P22"))
$ws.Range("P22").Value2 = 31.4

$ws.Range("O23").Copy($ws.Range("P23"))
$ws.Range("P23").Value2 = 31.5

$ws.Range("O24").Copy($ws.Range("P24"))
$ws.Range("P24").Value2 = 23.8

# --- Row 25 (bottom border row) ---
$ws.Range("O25").Copy($ws.Range("P25"))
$ws.Range("P25").Value2 = "-"

# --- Update selection to match target workbook state ---
$ws.Range("R9").Select()
